# Update the per-unit bus voltage magnitude results (vm_pu.xlsx) for the
# "case with 380 kV done" run: the slack-bus voltage setpoint in column B
# moved from 1.05 p.u. to 1.02 p.u., and the resulting bus voltages for
# all other buses (columns C-F, I-N) for every time step (rows 2-25) were
# recomputed accordingly. Column A (time index) and column G (unchanged at
# 1) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020113072211601
$ws.Range("D2").Value = 1.025907549841363
$ws.Range("E2").Value = 1.021144118774766
$ws.Range("F2").Value = 1.0329127437976
$ws.Range("I2").Value = 1.029301764090588
$ws.Range("J2").Value = 1.025312019603458
$ws.Range("K2").Value = 1.028731747364121
$ws.Range("L2").Value = 1.023982317307763
$ws.Range("M2").Value = 1.035716606163783
$ws.Range("N2").Value = 1.012444161774373

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021261962373333
$ws.Range("D3").Value = 1.026774506353881
$ws.Range("E3").Value = 1.022125041220126
$ws.Range("F3").Value = 1.035119657753549
$ws.Range("I3").Value = 1.029569287689169
$ws.Range("J3").Value = 1.026096915882789
$ws.Range("K3").Value = 1.029405944572635
$ws.Range("L3").Value = 1.024769151726096
$ws.Range("M3").Value = 1.03772865977159
$ws.Range("N3").Value = 1.01271042521562

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022004970882829
$ws.Range("D4").Value = 1.02733483879364
$ws.Range("E4").Value = 1.022759758497706
$ws.Range("F4").Value = 1.036541670074275
$ws.Range("I4").Value = 1.029740427995969
$ws.Range("J4").Value = 1.026603908144716
$ws.Range("K4").Value = 1.029840916129397
$ws.Range("L4").Value = 1.025277674085076
$ws.Range("M4").Value = 1.039024265715438
$ws.Range("N4").Value = 1.012882245700216

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022317238422022
$ws.Range("D5").Value = 1.027570249420303
$ws.Range("E5").Value = 1.023026593870859
$ws.Range("F5").Value = 1.037138079406133
$ws.Range("I5").Value = 1.029811906847539
$ws.Range("J5").Value = 1.026816836778333
$ws.Range("K5").Value = 1.030023473656696
$ws.Range("L5").Value = 1.025491311806165
$ws.Range("M5").Value = 1.039567455199408
$ws.Range("N5").Value = 1.012954367034769

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022369664122958
$ws.Range("D6").Value = 1.027609766967123
$ws.Range("E6").Value = 1.023071396759593
$ws.Range("F6").Value = 1.037238137761716
$ws.Range("I6").Value = 1.029823881043001
$ws.Range("J6").Value = 1.026852576113832
$ws.Range("K6").Value = 1.030054108061013
$ws.Range("L6").Value = 1.025527174054265
$ws.Range("M6").Value = 1.039658573032217
$ws.Range("N6").Value = 1.012966469969992

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022009143779442
$ws.Range("D7").Value = 1.027337984961147
$ws.Range("E7").Value = 1.022763323961796
$ws.Range("F7").Value = 1.036549644798842
$ws.Range("I7").Value = 1.029741384938109
$ws.Range("J7").Value = 1.026606754133681
$ws.Range("K7").Value = 1.029843356665969
$ws.Range("L7").Value = 1.025280529290542
$ws.Range("M7").Value = 1.039031529625752
$ws.Range("N7").Value = 1.012883209828348

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020501429456914
$ws.Range("D8").Value = 1.026200676677375
$ws.Range("E8").Value = 1.02147562805722
$ws.Range("F8").Value = 1.033659845311436
$ws.Range("I8").Value = 1.029392583080773
$ws.Range("J8").Value = 1.025577463898875
$ws.Range("K8").Value = 1.028959861354714
$ws.Range("L8").Value = 1.024248359316408
$ws.Range("M8").Value = 1.036397918301791
$ws.Range("N8").Value = 1.012534244302435

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.017841461747346
$ws.Range("D9").Value = 1.024191569295324
$ws.Range("E9").Value = 1.019206437175363
$ws.Range("F9").Value = 1.028520098064401
$ws.Range("I9").Value = 1.028762807382238
$ws.Range("J9").Value = 1.023756843387069
$ws.Range("K9").Value = 1.027393155065546
$ws.Range("L9").Value = 1.022424788849328
$ws.Range("M9").Value = 1.031707270879455
$ws.Range("N9").Value = 1.011915703180959

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016065843026293
$ws.Range("D10").Value = 1.022848677895947
$ws.Range("E10").Value = 1.01769347460794
$ws.Range("F10").Value = 1.025059466488303
$ws.Range("I10").Value = 1.028332647600177
$ws.Range("J10").Value = 1.022538358834727
$ws.Range("K10").Value = 1.026341932183038
$ws.Range("L10").Value = 1.021205781171598
$ws.Range("M10").Value = 1.028544615183169
$ws.Range("N10").Value = 1.011500877682333

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.015296395024898
$ws.Range("D11").Value = 1.022266338823768
$ws.Range("E11").Value = 1.017038282302081
$ws.Range("F11").Value = 1.023552391969258
$ws.Range("I11").Value = 1.028143909028165
$ws.Range("J11").Value = 1.022009591893838
$ws.Range("K11").Value = 1.025885113365466
$ws.Range("L11").Value = 1.020677133480675
$ws.Range("M11").Value = 1.027166262695799
$ws.Range("N11").Value = 1.011320661650387

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.015010495209666
$ws.Range("D12").Value = 1.022049900914824
$ws.Range("E12").Value = 1.016794902181066
$ws.Range("F12").Value = 1.022991263316605
$ws.Range("I12").Value = 1.028073428425668
$ws.Range("J12").Value = 1.021813008418751
$ws.Range("K12").Value = 1.025715182896988
$ws.Range("L12").Value = 1.020480646544326
$ws.Range("M12").Value = 1.02665290438406
$ws.Range("N12").Value = 1.011253631487985

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.015071825934772
$ws.Range("D13").Value = 1.02209633355786
$ws.Range("E13").Value = 1.016847108643818
$ws.Range("F13").Value = 1.023111688373933
$ws.Range("I13").Value = 1.028088563752399
$ws.Range("J13").Value = 1.02185518423901
$ws.Range("K13").Value = 1.025751644803217
$ws.Range("L13").Value = 1.020522799287751
$ws.Range("M13").Value = 1.026763084429602
$ws.Range("N13").Value = 1.011268013764834

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.015272764365266
$ws.Range("D14").Value = 1.022248450689658
$ws.Range("E14").Value = 1.017018164693048
$ws.Range("F14").Value = 1.023506036361868
$ws.Range("I14").Value = 1.028138090738538
$ws.Range("J14").Value = 1.021993345844298
$ws.Range("K14").Value = 1.025871071925191
$ws.Range("L14").Value = 1.020660894343568
$ws.Range("M14").Value = 1.027123856655217
$ws.Range("N14").Value = 1.011315122754274

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.015396556794967
$ws.Range("D15").Value = 1.022342157613433
$ws.Range("E15").Value = 1.017123556247852
$ws.Range("F15").Value = 1.02374882911644
$ws.Range("I15").Value = 1.028168556219573
$ws.Range("J15").Value = 1.022078448391129
$ws.Range("K15").Value = 1.02594462203283
$ws.Range("L15").Value = 1.020745962818837
$ws.Range("M15").Value = 1.027345956695186
$ws.Range("N15").Value = 1.011344136223606

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016116895963296
$ws.Range("D16").Value = 1.022887307544253
$ws.Range("E16").Value = 1.01773695592182
$ws.Range("F16").Value = 1.025159301313764
$ws.Range("I16").Value = 1.028345121157085
$ws.Range("J16").Value = 1.022573426818955
$ws.Range("K16").Value = 1.026372215149591
$ws.Range("L16").Value = 1.02124084855173
$ws.Range("M16").Value = 1.028635900998689
$ws.Range("N16").Value = 1.011512825456937

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.016568584426529
$ws.Range("D17").Value = 1.02322903430362
$ws.Range("E17").Value = 1.018121705258734
$ws.Range("F17").Value = 1.026041721655259
$ws.Range("I17").Value = 1.028455210801154
$ws.Range("J17").Value = 1.022883602811344
$ws.Range("K17").Value = 1.026639994319591
$ws.Range("L17").Value = 1.02155105917879
$ws.Range("M17").Value = 1.029442638115486
$ws.Range("N17").Value = 1.011618480285836

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.016831989647009
$ws.Range("D18").Value = 1.02342827509306
$ws.Range("E18").Value = 1.018346116472372
$ws.Range("F18").Value = 1.026555595802164
$ws.Range("I18").Value = 1.028519185456894
$ws.Range("J18").Value = 1.02306441193229
$ws.Range("K18").Value = 1.026796028042214
$ws.Range("L18").Value = 1.021731921842835
$ws.Range("M18").Value = 1.029912338041099
$ws.Range("N18").Value = 1.011680049709021

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.016921794451323
$ws.Range("D19").Value = 1.023496197103672
$ws.Range("E19").Value = 1.018422633845311
$ws.Range("F19").Value = 1.026730674651112
$ws.Range("I19").Value = 1.028540958739923
$ws.Range("J19").Value = 1.023126044363887
$ws.Range("K19").Value = 1.026849204893157
$ws.Range("L19").Value = 1.021793578177584
$ws.Range("M19").Value = 1.03007234953702
$ws.Range("N19").Value = 1.01170103360083

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.016520128451699
$ws.Range("D20").Value = 1.023192378821826
$ws.Range("E20").Value = 1.018080425990343
$ws.Range("F20").Value = 1.025947132118013
$ws.Range("I20").Value = 1.028443423940326
$ws.Range("J20").Value = 1.022850335390929
$ws.Range("K20").Value = 1.026611280440189
$ws.Range("L20").Value = 1.02151778459162
$ws.Range("M20").Value = 1.029356171700485
$ws.Range("N20").Value = 1.011607150449321

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.015213595588735
$ws.Range("D21").Value = 1.022203659627649
$ws.Range("E21").Value = 1.016967793312679
$ws.Range("F21").Value = 1.023389947878611
$ws.Range("I21").Value = 1.028123516640202
$ws.Range("J21").Value = 1.021952665569619
$ws.Range("K21").Value = 1.025835910469215
$ws.Range("L21").Value = 1.020620232221461
$ws.Range("M21").Value = 1.027017656554228
$ws.Range("N21").Value = 1.011301252818398

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014391589228909
$ws.Range("D22").Value = 1.021581252084589
$ws.Range("E22").Value = 1.016268163458277
$ws.Range("F22").Value = 1.021774404613471
$ws.Range("I22").Value = 1.027920209029932
$ws.Range("J22").Value = 1.021387246087505
$ws.Range("K22").Value = 1.025346970734335
$ws.Range("L22").Value = 1.02005518900257
$ws.Range("M22").Value = 1.025539353315115
$ws.Range("N22").Value = 1.011108402331057

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01482740243391
$ws.Range("D23").Value = 1.021911274972762
$ws.Range("E23").Value = 1.016639058154675
$ws.Range("F23").Value = 1.022631582648994
$ws.Range("I23").Value = 1.028028192742783
$ws.Range("J23").Value = 1.021687083021345
$ws.Range("K23").Value = 1.025606303655856
$ws.Range("L23").Value = 1.02035479770663
$ws.Range("M23").Value = 1.026323800163769
$ws.Range("N23").Value = 1.01121068559295

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.016542023797562
$ws.Range("D24").Value = 1.023208942111959
$ws.Range("E24").Value = 1.018099078335575
$ws.Range("F24").Value = 1.025989875611128
$ws.Range("I24").Value = 1.028448750653203
$ws.Range("J24").Value = 1.022865367850418
$ws.Range("K24").Value = 1.02662425549274
$ws.Range("L24").Value = 1.021532820186447
$ws.Range("M24").Value = 1.029395244795513
$ws.Range("N24").Value = 1.011612270091779

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.018529521387344
$ws.Range("D25").Value = 1.024711578161196
$ws.Range("E25").Value = 1.019793099892316
$ws.Range("F25").Value = 1.029854701415376
$ws.Range("I25").Value = 1.028927427011269
$ws.Range("J25").Value = 1.024228344119562
$ws.Range("K25").Value = 1.027799367876166
$ws.Range("L25").Value = 1.022896798905173
$ws.Range("M25").Value = 1.032926033057942
$ws.Range("N25").Value = 1.012076042975463
